$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry a leading "Test Number"/"1" column (column A) in
# front of the real Username/PassWord credential columns. Registration no
# longer needs that column, so drop it - this shifts Username/PassWord
# (formerly B/C) left into A/B.
$ws.Range("A1:A2").EntireColumn.Delete()

# Row 3 only held a stray leftover cell (A3) below the data - remove it too,
# so the sheet is a clean 2x2 Username/PassWord table.
$ws.Range("A3").EntireRow.Delete()

# Leave the selection where the author's last click landed.
$ws.Range("B7").Select()
